$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reverse_Leg_Lift")
Write-Host $ws.Name
